$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.761.11"
$ws.Range("E2").Value = "  +2.99%  "
$ws.Range("D3").Value = "2.697.44"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'526.01"
$ws.Range("E5").Value = "  +0.77%  "
$ws.Range("D6").Value = "'145.05"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "2.720.00"
$ws.Range("E9").Value = "  +2.68%  "
$ws.Range("D10").Value = "'6.70"
$ws.Range("E10").Value = "  +5.51%  "
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D14").Value = "3.194.77"
$ws.Range("E14").Value = "  +3.09%  "
$ws.Range("D15").Value = "60.729.72"
$ws.Range("E15").Value = "  +2.96%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'21.25"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.721.03"
$ws.Range("E17").Value = "  +2.94%  "
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'345.78"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'4.50"
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("D21").Value = "'10.60"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = "  +4.18%  "
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").Value = "'63.66"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  +3.47%  "
$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "0.0₃0817"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").Value = "'6.83"
$ws.Range("E30").Value = "  +9.02%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("D33").Value = "'19.04"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'150.16"
$ws.Range("E34").Value = "  +0.18%  "
$ws.Range("E35").Value = "  +6.12%  "
$ws.Range("E36").Value = "  +8.21%  "
$ws.Range("E37").Value = "  -3.69%  "
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").Value = "'37.15"
$ws.Range("E40").Value = "  +1.10%  "
$ws.Range("D41").Value = "'3.66"
$ws.Range("E41").Value = "  -0.58%  "
$ws.Range("D42").Value = "'282.22"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("D43").Value = "'20.08"
$ws.Range("E43").Value = "  +2.30%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "2.143.52"
$ws.Range("E45").Value = "  +7.76%  "
$ws.Range("D46").Value = "'0.0985"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'10.48"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  +0.92%  "

# Strip the quote-prefix marker so cells end up with no explicit style,
# matching the original (unstyled) data cells.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D50").ClearFormats()
